$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:D3").Value = -0.0371
$ws.Range("G2:G3").Value = 0.008403575989782886
$ws.Range("H2:H3").Value = 0.007356321839080459
$ws.Range("I2:I3").Value = -0.05261813537675607
$ws.Range("J2:J3").Value = -0.05261813537675607
$ws.Range("K2:K3").Value = -3.77
$ws.Range("L2:L3").Value = -0.04814814814814815

# Note: O2/R2 go from 0 -> -0 and O3/R3 go from -0 -> 0 (sign-of-zero only;
# numerically identical to 0). Pass "-0" as a string so the literal isn't
# collapsed to plain 0 before assignment.
$ws.Range("O2").Value = "-0"
$ws.Range("O3").Value = 0
$ws.Range("R2").Value = "-0"
$ws.Range("R3").Value = 0

$ws.Range("U2:U3").Value = 3.48
$ws.Range("V2:V3").Value = 0.04793388429752066
$ws.Range("W2:W3").Value = -0.06744186046511629
$ws.Range("X2:X3").Value = 0.07657097999753487
$ws.Range("Y2:Y3").Value = -0.1440128404626512
$ws.Range("Z2:Z3").Value = 0.9446254071661238
$ws.Range("AA2:AA3").Value = -0.04970442755459042
$ws.Range("AB2:AB3").Value = 0.05763883659941507
$ws.Range("AC2:AC3").Value = -0.1073432641540055
$ws.Range("AD2:AD3").Value = 42
$ws.Range("AF2:AF3").Value = 42
$ws.Range("AG2:AG3").Value = 38.52
$ws.Range("AH2:AH3").Value = 0.3664921465968586
$ws.Range("AI2:AI3").Value = 0.4276985743380855
$ws.Range("AJ2:AJ3").Value = 0.3466522678185746
$ws.Range("AK2:AK3").Value = 0.4066722972972973
$ws.Range("AL2:AL3").Value = 0.505
$ws.Range("AM2:AM3").Value = 0.374
$ws.Range("AN2:AN3").Value = -17.87234042553191
$ws.Range("AO2:AO3").Value = -8.158415841584159
$ws.Range("AP2:AP3").Value = -16.39148936170213
$ws.Range("AQ2:AQ3").Value = -11.01604278074866
